# Auto-generated Word COM-interop script to apply the resume content update
$d = $word.ActiveDocument

# --- Step 1: Update existing paragraph text in place (1-based Paragraphs index) ---
$d.Paragraphs.Item(4).Range.Text = 'Distinguished Research and Data Analytics Leader with 21 years of experience directing groundbreaking applied research projects that have shaped policy, influenced elections, and transformed community development outcomes. Led multi-million dollar research initiatives serving thousands of analysts nationwide, with proven expertise in translating complex research insights for diverse stakeholders including elected officials, government agencies, and community organizations. Expert in research methodology design, statistical analysis, and community partnership development with extensive experience addressing systemic poverty challenges and delivering evidence-based solutions that drive meaningful social impact.'
$d.Paragraphs.Item(10).Range.Text = 'PARTNER - Siege Analytics, Washington, DC | 2005 – Present'
$d.Paragraphs.Item(12).Range.Text = '• Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations'
$d.Paragraphs.Item(13).Range.Text = '• Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics'
$d.Paragraphs.Item(14).Range.Text = '• Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets'
$d.Paragraphs.Item(15).Range.Text = '• Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering'
$d.Paragraphs.Item(16).Range.Text = '• Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications'
$d.Paragraphs.Item(17).Range.Text = '• Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices'
$d.Paragraphs.Item(18).Range.Text = 'DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | June 2021 – May 2023'
$d.Paragraphs.Item(19).Range.Text = 'Data Platform Development and Team Leadership'
$d.Paragraphs.Item(20).Range.Text = '• Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES'
$d.Paragraphs.Item(21).Range.Text = '• Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions'
$d.Paragraphs.Item(22).Range.Text = '• Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI'
$d.Paragraphs.Item(23).Range.Text = '• Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company''s distinguishing products'
$d.Paragraphs.Item(24).Range.Text = 'SOFTWARE ENGINEER - Mautinoa Technologies, Washington, DC | August 2016 – February 2018'
$d.Paragraphs.Item(25).Range.Text = 'Financial Technology and Humanitarian Crisis Solutions'
$d.Paragraphs.Item(26).Range.Text = '• Developed SimCrisis, a GeoDjango web application using Python, PostgreSQL/PostGIS, and NetLogo for multi-agent modeling and econometric simulations of crisis economies'
$d.Paragraphs.Item(27).Range.Text = '• Built modular application using Python, Django, and GRASS accepting rules extensions for ethnic strife, different crisis types, supply failures, and disaster scenarios'
$d.Paragraphs.Item(28).Range.Text = '• Liaised with officers from International Federation of Red Cross, UNICEF, and Chaos Communications Congress to improve platform using Docker and Ubuntu'
$d.Paragraphs.Item(29).Range.Text = '• Conceived and built application using Python, Pandas, and Jupyter to predict how crisis economies respond to different humanitarian interventions'
$d.Paragraphs.Item(30).Range.Text = 'SENIOR ANALYST - Myers Research, Washington, DC | August 2012 – February 2014'
$d.Paragraphs.Item(31).Range.Text = 'Quantitative and Qualitative Research for Democratic Campaigns'
$d.Paragraphs.Item(32).Range.Text = '• Developed RACSO, a web application for pollsters to fully administer research including questionnaire creation, versioning, and reporting'
$d.Paragraphs.Item(33).Range.Text = '• Wrote RFP and analyzed bids from 1,200 vendors before selecting implementation partner'
$d.Paragraphs.Item(34).Range.Text = '• Built prototype in R for comprehensive polling administration and sample file management'
$d.Paragraphs.Item(35).Range.Text = '• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research'
$d.Paragraphs.Item(36).Range.Text = 'RESEARCH DIRECTOR - Progressive Change Campaign Committee, Washington, DC | August 2011 – August 2012'
$d.Paragraphs.Item(42).Range.Text = 'SOFTWARE ENGINEER - Salsa Labs, Inc., Washington, DC | January 2011 – August 2011'
$d.Paragraphs.Item(43).Range.Text = 'Political Technology Development'
$d.Paragraphs.Item(44).Range.Text = '• Developed software solutions for political campaigns and advocacy groups'
$d.Paragraphs.Item(45).Range.Text = '• Built web applications for voter engagement and campaign management'
$d.Paragraphs.Item(46).Range.Text = '• Integrated third-party APIs and data sources for campaign tools'
$d.Paragraphs.Item(47).Range.Text = '• Collaborated with political strategists to translate requirements into technical solutions'
$d.Paragraphs.Item(48).Range.Text = 'INTERIM TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | April 2009 – October 2009'
$d.Paragraphs.Item(49).Range.Text = 'Nonprofit Technology Integration'
$d.Paragraphs.Item(50).Range.Text = '• Integrated technology solutions within organizational frameworks for social justice organizations'
$d.Paragraphs.Item(51).Range.Text = '• Developed data management systems for community organizing efforts'
$d.Paragraphs.Item(52).Range.Text = '• Provided technical training and support to nonprofit staff'
$d.Paragraphs.Item(53).Range.Text = '• Built custom applications for community engagement and advocacy'
$d.Paragraphs.Item(54).Range.Text = 'PROGRAMMER - Lake Research Partners, Washington, DC | April 2008 – December 2008'
$d.Paragraphs.Item(56).Range.Text = '• Developed data analysis tools for political polling and research'
$d.Paragraphs.Item(57).Range.Text = '• Built statistical models for voter behavior analysis'
$d.Paragraphs.Item(58).Range.Text = '• Created data visualization tools for research presentations'
$d.Paragraphs.Item(59).Range.Text = '• Supported senior researchers with technical analysis and reporting'
$d.Paragraphs.Item(60).Range.Text = 'FIELD DIRECTOR - The Feldman Group, Washington, DC | August 2007 – April 2008'
$d.Paragraphs.Item(62).Range.Text = '• Managed field operations for political campaigns and research projects'
$d.Paragraphs.Item(63).Range.Text = '• Developed data collection and management systems for field work'
$d.Paragraphs.Item(64).Range.Text = '• Trained field staff on data collection protocols and quality control'
$d.Paragraphs.Item(65).Range.Text = '• Analyzed field data to inform campaign strategy and research findings'
$d.Paragraphs.Item(67).Range.Text = 'Software Development and Innovation'
$d.Paragraphs.Item(68).Range.Text = '• Conceived and deployed redistricting software used by thousands of analysts nationwide'
$d.Paragraphs.Item(69).Range.Text = '• Developed boundary estimation system using incomplete data without ML requirements'
$d.Paragraphs.Item(70).Range.Text = '• Created econometric simulation platform for humanitarian intervention modeling'

# --- Step 2: Insert new bullet after the 'five-year strategic plans' bullet (was index 23) ---
$p23 = $d.Paragraphs.Item(23)
$p23.Range.InsertParagraphAfter()
$d.Paragraphs.Item(24).Range.Text = '• Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices'

# --- Step 3: Append new bullet at the very end of the document (after last paragraph) ---
$lastIdx = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIdx)
$pLast.Range.InsertParagraphAfter()
$d.Paragraphs.Item($lastIdx + 1).Range.Text = '• Built comprehensive survey operations platform from RFP through deployment'

Write-Host 'Edit complete'
